# Generate Report for handoff
#
# A new localization entry (a8d7baa1-8806-428c-92bd-f19d05061f94.md) shows up with
# status "Handoff transform failed" ahead of the existing ".localization-config"
# row, on all three sheets (Overview / zh-cn / de-de). The previously "Ready for
# handoff" entry also rolled forward to a new source uuid / handoff package /
# timestamps (eb9f0baa... -> 735e2a3e...).

$wb = $excel.ActiveWorkbook

$oldUuid  = "eb9f0baa-c02b-4d68-aef1-e20fe61f8e97"
$newUuid  = "735e2a3e-aa66-4b58-9612-6255bd4eaca9"
$failUuid = "a8d7baa1-8806-428c-92bd-f19d05061f94"
$oldHash  = "c6beba2b279d34f7f50a60d41ef92abac9194cf0"
$newHash  = "cab1c72d2cf683739f2b8d96785b874c26c91cb3"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/a88ffa8a7d9f5a3cc80751bf2fa3581f7cfab42e/e2e/"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a88ffa8a7d9f5a3cc80751bf2fa3581f7cfab42e/.localization-config"

$hyperlinkColor = 15570276 # OLE BGR for RGB(0x64,0x95,0xED) -> matches the workbook's existing custom HyperLink font

function Style-AsHyperlink($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = 1
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Hyperlinks.Delete() on any cell's Hyperlinks collection clears every
# hyperlink on the sheet here, so just nuke them all up front and re-add.
$ov.Range("A1").Hyperlinks.Delete()

$ov.Range("A2").Value = "$newUuid.md"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"

$ov.Range("A3").Value = "$failUuid.md"
$ov.Range("B3").Value = "Handoff transform failed"
$ov.Range("C3").Value = "Handoff transform failed"

$ov.Range("A4").Value = ".localization-config"
$ov.Range("B4").Value = "Not to be localized"
$ov.Range("C4").Value = "Not to be localized"

$ov.Hyperlinks.Add($ov.Range("A2"), "$mdBase$newUuid.md", "", "", "$newUuid.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "$mdBase$failUuid.md", "", "", "$failUuid.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

Style-AsHyperlink $ov.Range("A2")
Style-AsHyperlink $ov.Range("A3")
Style-AsHyperlink $ov.Range("A4")

# ---------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A1").Hyperlinks.Delete()

$zh.Range("A2").Value = "$newUuid.md"
$zh.Range("B2").Value = "Ready for handoff"
$zh.Range("C2").Value = "$newUuid.$newHash.zh-cn.xlf"
$zh.Range("D2").Value = "2016-02-19 06:04:09"
$zh.Range("G2").Value = "0001-01-01 00:00:00"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "$failUuid.md"
$zh.Range("B3").Value = "Handoff transform failed"
$zh.Range("D3").Value = "0001-01-01 00:00:00"
$zh.Range("G3").Value = "0001-01-01 00:00:00"
$zh.Range("H3").Value = "Ignored"

$zh.Range("A4").Value = ".localization-config"
$zh.Range("B4").Value = "Not to be localized"
$zh.Range("D4").Value = "0001-01-01 00:00:00"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Ignored"

$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bb45c1ad3a92060fef75a7df397d02a5449e5832/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/$newUuid.$newHash.zh-cn.xlf"

$zh.Hyperlinks.Add($zh.Range("A2"), "$mdBase$newUuid.md", "", "", "$newUuid.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), $zhXlfUrl, "", "", "$newUuid.$newHash.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "$mdBase$failUuid.md", "", "", "$failUuid.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

Style-AsHyperlink $zh.Range("A2")
Style-AsHyperlink $zh.Range("C2")
Style-AsHyperlink $zh.Range("A3")
Style-AsHyperlink $zh.Range("A4")

$zh.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("A1").Hyperlinks.Delete()

$de.Range("A2").Value = "$newUuid.md"
$de.Range("B2").Value = "Ready for handoff"
$de.Range("C2").Value = "$newUuid.$newHash.de-de.xlf"
$de.Range("D2").Value = "2016-02-19 06:04:22"
$de.Range("G2").Value = "0001-01-01 00:00:00"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "$failUuid.md"
$de.Range("B3").Value = "Handoff transform failed"
$de.Range("D3").Value = "0001-01-01 00:00:00"
$de.Range("G3").Value = "0001-01-01 00:00:00"
$de.Range("H3").Value = "Ignored"

$de.Range("A4").Value = ".localization-config"
$de.Range("B4").Value = "Not to be localized"
$de.Range("D4").Value = "0001-01-01 00:00:00"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Ignored"

$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b1bab51a9ed686bb53b1ef0cf3753768d24f8a52/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/$newUuid.$newHash.de-de.xlf"

$de.Hyperlinks.Add($de.Range("A2"), "$mdBase$newUuid.md", "", "", "$newUuid.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), $deXlfUrl, "", "", "$newUuid.$newHash.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "$mdBase$failUuid.md", "", "", "$failUuid.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

Style-AsHyperlink $de.Range("A2")
Style-AsHyperlink $de.Range("C2")
Style-AsHyperlink $de.Range("A3")
Style-AsHyperlink $de.Range("A4")

$de.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
